$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the bibliographic citation column (N) with the full citation text
# for each dataset row (previously these cells held stray column-name
# references, or were left empty).
$ws.Range("N2").Value = "Seebens, H. et al. No saturation in the accumulation of alien species worldwide. Nat. Commun. 8, 14435 (2017)"
$ws.Range("N3").Value = "Pagad, S. et al. Country Compendium of the Global Register of Introduced and Invasive Species. Sci. Data 9, 391 (2022)"
$ws.Range("N4").Value = "Capinha, C. et al. Diversity, biogeography and the global flows of alien amphibians and reptiles. Divers. Distrib. 23, 1313–1322 (2017)"
$ws.Range("N5").Value = "van Kleunen, M. et al. The Global Naturalized Alien Flora (GloNAF) database. Ecology 100, e02542 (2019)"
$ws.Range("N6").Value = "Dyer, E. E., Redding, D. W. & Blackburn, T. M. The global avian invasions atlas, a database of alien bird distributions worldwide. Sci. Data 4, 170041 (2017)"
$ws.Range("N7").Value = "IUCN. The IUCN Red List of Threatened Species. https://www.iucnredlist.org (2025)"

# Re-select/scroll the frozen pane view to match the edited area.
$ws.Range("A1").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("F1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("N7").Select()

Write-Output "done"
